# Add working set of sequences
# Updates the cue sheet (Sheet1) so that each row's count (column B),
# image (column C), word (column D) and category (column E) reflect the
# new working set of sequences.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Count = 116; Image = "flower/flower019.jpg"; Word = "pflegen";  Category = "flower" },
    @{ Row = 3;  Count = 91;  Image = "flower/flower013.jpg"; Word = "posten";   Category = "flower" },
    @{ Row = 4;  Count = 83;  Image = "flower/flower009.jpg"; Word = "formen";   Category = "flower" },
    @{ Row = 5;  Count = 109; Image = "flower/flower001.jpg"; Word = "schicken"; Category = "flower" },
    @{ Row = 6;  Count = 71;  Image = "dog/dog020.jpg";       Word = "drehen";   Category = "dog" },
    @{ Row = 7;  Count = 59;  Image = "flower/flower014.jpg"; Word = "tauschen"; Category = "flower" },
    @{ Row = 8;  Count = 60;  Image = "flower/flower006.jpg"; Word = "tagen";    Category = "flower" },
    @{ Row = 9;  Count = 81;  Image = "flower/flower016.jpg"; Word = "schätzen"; Category = "flower" },
    @{ Row = 10; Count = 32;  Image = "flower/flower008.jpg"; Word = "antun";    Category = "flower" },
    @{ Row = 11; Count = 126; Image = "flower/flower010.jpg"; Word = "husten";   Category = "flower" },
    @{ Row = 12; Count = 86;  Image = "flower/flower017.jpg"; Word = "klappen";  Category = "flower" },
    @{ Row = 13; Count = 52;  Image = "flower/flower031.jpg"; Word = "biegen";   Category = "flower" },
    @{ Row = 14; Count = 19;  Image = "dog/dog004.jpg";       Word = "fühlen";   Category = "dog" },
    @{ Row = 15; Count = 41;  Image = "dog/dog005.jpg";       Word = "nehmen";   Category = "dog" },
    @{ Row = 16; Count = 50;  Image = "flower/flower018.jpg"; Word = "krachen";  Category = "flower" },
    @{ Row = 17; Count = 97;  Image = "dog/dog019.jpg";       Word = "rücken";   Category = "dog" },
    @{ Row = 18; Count = 65;  Image = "dog/dog008.jpg";       Word = "wenden";   Category = "dog" },
    @{ Row = 19; Count = 100; Image = "dog/dog027.jpg";       Word = "enden";    Category = "dog" },
    @{ Row = 20; Count = 6;   Image = "dog/dog030.jpg";       Word = "runden";   Category = "dog" },
    @{ Row = 21; Count = 24;  Image = "dog/dog003.jpg";       Word = "spielen";  Category = "dog" },
    @{ Row = 22; Count = 90;  Image = "dog/dog026.jpg";       Word = "rasen";    Category = "dog" },
    @{ Row = 23; Count = 79;  Image = "dog/dog023.jpg";       Word = "hoffen";   Category = "dog" },
    @{ Row = 24; Count = 30;  Image = "dog/dog014.jpg";       Word = "sieben";   Category = "dog" },
    @{ Row = 25; Count = 111; Image = "flower/flower027.jpg"; Word = "loben";    Category = "flower" },
    @{ Row = 26; Count = 92;  Image = "dog/dog002.jpg";       Word = "regnen";   Category = "dog" },
    @{ Row = 27; Count = 122; Image = "flower/flower030.jpg"; Word = "wiegen";   Category = "flower" },
    @{ Row = 28; Count = 117; Image = "dog/dog029.jpg";       Word = "haken";    Category = "dog" },
    @{ Row = 29; Count = 63;  Image = "flower/flower020.jpg"; Word = "gelten";   Category = "flower" },
    @{ Row = 30; Count = 0;   Image = "dog/dog016.jpg";       Word = "kaufen";   Category = "dog" },
    @{ Row = 31; Count = 123; Image = "dog/dog001.jpg";       Word = "liefern";  Category = "dog" },
    @{ Row = 32; Count = 5;   Image = "flower/flower007.jpg"; Word = "fliegen";  Category = "flower" },
    @{ Row = 33; Count = 95;  Image = "dog/dog010.jpg";       Word = "langen";   Category = "dog" }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.Count
    $ws.Range("C" + $r.Row).Value = $r.Image
    $ws.Range("D" + $r.Row).Value = $r.Word
    $ws.Range("E" + $r.Row).Value = $r.Category
}
